# Update scraped Market Board price/profit figures across all Job sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the latest scheduled price pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 481.5
$ws.Range("I28").Value = 337.8
$ws.Range("J28").Value = 1200
$ws.Range("K28").Value = 337.8
$ws.Range("L28").Value = 1200
$ws.Range("M28").Value = 147.2
$ws.Range("N28").Value = -2170

$ws.Range("H100").Value = 3663.25
$ws.Range("I100").Value = 800
$ws.Range("J100").Value = 4072.2856
$ws.Range("K100").Value = 800
$ws.Range("L100").Value = 4072.2856
$ws.Range("M100").Value = -259
$ws.Range("N100").Value = -5154.2856

$ws.Range("H111").Value = 2355.3333
$ws.Range("I111").Value = 3183
$ws.Range("K111").Value = 9549
$ws.Range("M111").Value = -6482

$ws.Range("H116").Value = 4449.875
$ws.Range("I116").Value = 2313
$ws.Range("K116").Value = 2313
$ws.Range("M116").Value = 1129

$ws.Range("H129").Value = 270998.62
$ws.Range("J129").Value = 345675.94
$ws.Range("L129").Value = 1037027.82
$ws.Range("N129").Value = -1047027.82

$ws.Range("H132").Value = 3416.3076
$ws.Range("I132").Value = 3617.7083
$ws.Range("K132").Value = 10853.1249
$ws.Range("M132").Value = -8323.124899999999

$ws.Range("H137").Value = 1740
$ws.Range("I137").Value = 1322.3529
$ws.Range("J137").Value = 3160
$ws.Range("K137").Value = 3967.0587
$ws.Range("L137").Value = 9480
$ws.Range("M137").Value = -1417.0587
$ws.Range("N137").Value = -14580

$ws.Range("H138").Value = 1900.9495
$ws.Range("I138").Value = 1351.9642
$ws.Range("J138").Value = 2117.4507
$ws.Range("K138").Value = 4055.8926
$ws.Range("L138").Value = 6352.3521
$ws.Range("M138").Value = 1084.1074
$ws.Range("N138").Value = -16632.3521

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3214.84
$ws.Range("I32").Value = 2681.0652
$ws.Range("K32").Value = 2681.0652
$ws.Range("M32").Value = -2394.0652

$ws.Range("H74").Value = 100001560
$ws.Range("I74").Value = 166667600
$ws.Range("J74").Value = 2500
$ws.Range("K74").Value = 166667600
$ws.Range("L74").Value = 2500
$ws.Range("M74").Value = -166666726
$ws.Range("N74").Value = -4248

$ws.Range("H76").Value = 15000
$ws.Range("J76").Value = 15000
$ws.Range("L76").Value = 15000
$ws.Range("N76").Value = -15676

$ws.Range("H77").Value = 100001560
$ws.Range("I77").Value = 166667600
$ws.Range("J77").Value = 2500
$ws.Range("K77").Value = 833338000
$ws.Range("L77").Value = 12500
$ws.Range("M77").Value = -833333632
$ws.Range("N77").Value = -21236

$ws.Range("H79").Value = 15000
$ws.Range("J79").Value = 15000
$ws.Range("L79").Value = 15000
$ws.Range("N79").Value = -17340

$ws.Range("H102").Value = 2150
$ws.Range("I102").Value = 2150
$ws.Range("K102").Value = 2150
$ws.Range("M102").Value = -528

$ws.Range("H110").Value = 792.75
$ws.Range("I110").Value = 702.3
$ws.Range("J110").Value = 1245
$ws.Range("K110").Value = 702.3
$ws.Range("L110").Value = 1245
$ws.Range("M110").Value = 1342.7
$ws.Range("N110").Value = -5335

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 463.33334
$ws.Range("I22").Value = 463.33334
$ws.Range("K22").Value = 463.33334
$ws.Range("M22").Value = -290.33334

$ws.Range("H100").Value = 29571.5
$ws.Range("J100").Value = 29571.5
$ws.Range("L100").Value = 29571.5
$ws.Range("N100").Value = -31735.5

$ws.Range("H105").Value = 1854577.2
$ws.Range("I105").Value = 2551.7144
$ws.Range("K105").Value = 2551.7144
$ws.Range("M105").Value = -804.7143999999998

$ws.Range("H107").Value = 913.5625
$ws.Range("I107").Value = 864.1111
$ws.Range("J107").Value = 977.1429000000001
$ws.Range("K107").Value = 864.1111
$ws.Range("L107").Value = 977.1429000000001
$ws.Range("M107").Value = 1055.8889
$ws.Range("N107").Value = -4817.1429

$ws.Range("H134").Value = 3908.7
$ws.Range("I134").Value = 4449.32
$ws.Range("J134").Value = 1205.6
$ws.Range("K134").Value = 13347.96
$ws.Range("L134").Value = 3616.8
$ws.Range("M134").Value = -10812.96
$ws.Range("N134").Value = -8686.799999999999

$ws.Range("H138").Value = 50000
$ws.Range("J138").Value = 50000
$ws.Range("L138").Value = 50000
$ws.Range("N138").Value = -60280

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 25236.334
$ws.Range("I58").Value = 1728.1666
$ws.Range("J58").Value = 56580.555
$ws.Range("K58").Value = 1728.1666
$ws.Range("L58").Value = 56580.555
$ws.Range("M58").Value = -1525.1666
$ws.Range("N58").Value = -56986.555

$ws.Range("H86").Value = 25673
$ws.Range("I86").Value = 3993.3333
$ws.Range("J86").Value = 47352.668
$ws.Range("K86").Value = 3993.3333
$ws.Range("L86").Value = 47352.668
$ws.Range("M86").Value = -2870.3333
$ws.Range("N86").Value = -49598.668

$ws.Range("H89").Value = 25673
$ws.Range("I89").Value = 3993.3333
$ws.Range("J89").Value = 47352.668
$ws.Range("K89").Value = 19966.6665
$ws.Range("L89").Value = 236763.34
$ws.Range("M89").Value = -14350.6665
$ws.Range("N89").Value = -247995.34

$ws.Range("H105").Value = 3000
$ws.Range("I105").Value = 2000
$ws.Range("J105").Value = 3500
$ws.Range("K105").Value = 2000
$ws.Range("L105").Value = 3500
$ws.Range("M105").Value = -253
$ws.Range("N105").Value = -6994

$ws.Range("H107").Value = 642.5
$ws.Range("I107").Value = 270.83334
$ws.Range("J107").Value = 1014.1667
$ws.Range("K107").Value = 270.83334
$ws.Range("L107").Value = 1014.1667
$ws.Range("M107").Value = 1649.16666
$ws.Range("N107").Value = -4854.1667

$ws.Range("H134").Value = 1608.3334
$ws.Range("I134").Value = 1066.2307
$ws.Range("K134").Value = 3198.6921
$ws.Range("M134").Value = -663.6921000000002

$ws.Range("H136").Value = 25236.334
$ws.Range("I136").Value = 1728.1666
$ws.Range("J136").Value = 56580.555
$ws.Range("K136").Value = 5184.4998
$ws.Range("L136").Value = 169741.665
$ws.Range("M136").Value = -2634.4998
$ws.Range("N136").Value = -174841.665

$ws.Range("H138").Value = 50000
$ws.Range("J138").Value = 50000
$ws.Range("L138").Value = 50000
$ws.Range("N138").Value = -60280

$ws.Range("H140").Value = 50000
$ws.Range("J140").Value = 50000
$ws.Range("L140").Value = 50000
$ws.Range("N140").Value = -60360

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 106.44444
$ws.Range("I12").Value = 49.5
$ws.Range("J12").Value = 122.71429
$ws.Range("K12").Value = 148.5
$ws.Range("L12").Value = 368.14287
$ws.Range("M12").Value = 24.5
$ws.Range("N12").Value = -714.14287

$ws.Range("H38").Value = 82.40000000000001
$ws.Range("I38").Value = 82.40000000000001
$ws.Range("K38").Value = 247.2
$ws.Range("M38").Value = 99.79999999999998

$ws.Range("H117").Value = 3300
$ws.Range("I117").Value = 1000
$ws.Range("J117").Value = 3555.5557
$ws.Range("K117").Value = 3000
$ws.Range("L117").Value = 10666.6671
$ws.Range("M117").Value = 442
$ws.Range("N117").Value = -17550.6671

$ws.Range("H131").Value = 137788.83
$ws.Range("J131").Value = 143650.36
$ws.Range("L131").Value = 430951.08
$ws.Range("N131").Value = -441031.08

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H74").Value = 51533.332
$ws.Range("J74").Value = 51533.332
$ws.Range("L74").Value = 51533.332
$ws.Range("N74").Value = -53405.332

$ws.Range("H77").Value = 51533.332
$ws.Range("J77").Value = 51533.332
$ws.Range("L77").Value = 154599.996
$ws.Range("N77").Value = -163959.996

$ws.Range("H113").Value = 3026.6667
$ws.Range("I113").Value = 2440
$ws.Range("J113").Value = 4200
$ws.Range("K113").Value = 2440
$ws.Range("L113").Value = 4200
$ws.Range("M113").Value = -270
$ws.Range("N113").Value = -8540

$ws.Range("H122").Value = 2371.4285
$ws.Range("I122").Value = 2600
$ws.Range("J122").Value = 2066.6667
$ws.Range("K122").Value = 7800
$ws.Range("L122").Value = 6200.000100000001
$ws.Range("M122").Value = -5350
$ws.Range("N122").Value = -11100.0001

$ws.Range("H132").Value = 42301.57
$ws.Range("I132").Value = 7472.4
$ws.Range("K132").Value = 22417.2
$ws.Range("M132").Value = -19887.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").ClearContents()  # value no longer available; cell removed

$ws.Range("H40").Value = 5999.8335
$ws.Range("I40").Value = 3748.75
$ws.Range("J40").Value = 6643
$ws.Range("K40").Value = 3748.75
$ws.Range("L40").Value = 6643
$ws.Range("M40").Value = -3612.75
$ws.Range("N40").Value = -6915

$ws.Range("H122").Value = 1228458
$ws.Range("I122").Value = 1785058.5
$ws.Range("J122").Value = 3937
$ws.Range("K122").Value = 5355175.5
$ws.Range("L122").Value = 11811
$ws.Range("M122").Value = -5352725.5
$ws.Range("N122").Value = -16711

$ws.Range("H132").Value = 1667.2609
$ws.Range("I132").Value = 1049.9474
$ws.Range("K132").Value = 3149.8422
$ws.Range("M132").Value = -619.8422

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 703.0909
$ws.Range("I107").Value = 714.8889
$ws.Range("J107").Value = 650
$ws.Range("K107").Value = 2144.6667
$ws.Range("L107").Value = 1950
$ws.Range("M107").Value = -224.6667000000002
$ws.Range("N107").Value = -5790

$ws.Range("H113").Value = 1338
$ws.Range("I113").Value = 1570.8889
$ws.Range("J113").Value = 290
$ws.Range("K113").Value = 4712.6667
$ws.Range("L113").Value = 870
$ws.Range("M113").Value = -2542.6667
$ws.Range("N113").Value = -5210

$ws.Range("H132").Value = 1137.871
$ws.Range("I132").Value = 774.5714
$ws.Range("K132").Value = 2323.7142
$ws.Range("M132").Value = 206.2857999999997
